# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2404
#   *_new  -> *_FV2410
# then turn the sheet into a proper Excel table ("Table1") with an
# autofilter, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base names shared by both the "FV2404" (left/old) and "FV2410" (right/new)
# column blocks; column K ("diff") sits between the two blocks and keeps its
# name unchanged.
$baseNames = @(
  "Segmentname",
  "Segmentgruppe",
  "Segment",
  "Datenelement",
  "Segment ID",
  "Code",
  "Qualifier",
  "Beschreibung",
  "Bedingungsausdruck",
  "Bedingung"
)

$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
  $ws.Range($leftCols[$i]  + "1").Value = $baseNames[$i] + "_FV2404"
  $ws.Range($rightCols[$i] + "1").Value = $baseNames[$i] + "_FV2410"
}

# Wrap the data range in a native Excel Table with an AutoFilter so the
# header row names above become the table's column headers.
$tableRange = $ws.Range("A1:U64")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row (split below row 1, top-left cell of the scrolling
# pane is A2).
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
